$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Define a print area on Sheet1 (adds the workbook-level _xlnm.Print_Area defined name)
$ws.PageSetup.PrintArea = '$A$2:$P$37'

# Update the "Fit to" / scale percentage used when printing
$ws.PageSetup.Zoom = 71

# Reflect the new selection (matches the print area) as the sheet's active selection
$ws.Range("A2:P37").Select() | Out-Null
